# Apply crypto price/volume/listing updates (GitHub Actions "Updated cryptos list" run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.842.79'
$ws.Range('D3').Value = '1.630.59'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  +0.67%  '
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('E7').Value = '  +0.58%  '
$ws.Range('D8').Value = "'0.254"
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').Value = "'0.0631"
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('D10').Value = "'19.53"
$ws.Range('E10').Value = '  -0.23%  '
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').Value = '1.857.27'
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('D14').Value = '1.634.94'
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('E15').Value = '  -1.30%  '
$ws.Range('D16').Value = '0.0₃0753'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').Value = "'62.61"
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '25.859.79'
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('D20').Value = "'4.38"
$ws.Range('E20').Value = '  -0.92%  '
$ws.Range('D21').Value = "'193.08"
$ws.Range('E21').Value = '  +1.30%  '
$ws.Range('D22').Value = "'9.90"
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = "'6.23"
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('D25').Value = "'143.37"
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('E26').Value = '  +0.62%  '
$ws.Range('E27').Value = '  +2.70%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = "'15.41"
$ws.Range('E29').Value = '  -0.45%  '
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('E31').Value = '  +0.99%  '
$ws.Range('E32').Value = '  -0.68%  '
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('D34').Value = "'1.55"
$ws.Range('E34').Value = '  -2.22%  '
$ws.Range('E35').Value = '  +1.69%  '
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').Value = '1.138.11'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('E38').Value = '  +0.77%  '
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('E40').Value = '  +0.48%  '
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('D42').Value = "'98.97"
$ws.Range('E42').Value = '  -1.31%  '
$ws.Range('E43').Value = '  -2.67%  '
$ws.Range('D44').Value = "'0.793"
$ws.Range('E44').Value = '  -0.61%  '
$ws.Range('D45').Value = '1.766.92'
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'56.24"
$ws.Range('E46').Value = '  +1.54%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = "'0.0527"
$ws.Range('E47').Value = '  +3.01%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = "'1.45"
$ws.Range('E48').Value = '  -1.33%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = "'0.415"
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'7.64"
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = "'0.0959"
$ws.Range('E51').Value = '  +0.49%  '
